# daily auto push: 2026-01-07 13:49 UTC
# Insert a new daily-ranking record (2026/01/07, 水, 19, 24) as a new row
# at sheet row 575, pushing the existing rows 575:616 down to 576:617.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 575; Excel shifts rows
# 575:616 down to 576:617 (and the sheet dimension grows to D617).
$ws.Rows(575).Insert()

# Column A holds dates as plain text (e.g. "2026/12/29"), not real Excel
# date serials, in this workbook. Force the new cell to Text first so the
# "2026/01/07" string isn't auto-converted into a date value, then drop
# the explicit number-format again so the cell matches its neighbours.
$ws.Range("A575").NumberFormat = "@"
$ws.Range("A575").Value = "2026/01/07"
$ws.Range("A575").ClearFormats()

$ws.Range("B575").Value = "水"
$ws.Range("C575").Value = 19
$ws.Range("D575").Value = 24
